# tambahan kode festronik saat export import dokumen
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header "Kode Festronik" in E1, matching the style of the other header cells (D1)
$ws.Range("E1").Value = "Kode Festronik"

# Copy full formatting (font, borders, alignment) from D1 to E1 to mirror the existing headers
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Move the active selection like the saved workbook (cosmetic, matches diff)
$ws.Range("J8").Select()
